$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 270-274 hold the field_license_wbddh controlled-vocabulary
# mapping. Column B/C values are entered first (row by row) so the shared
# string table picks up the same ordering as in the authored workbook; the
# A-column label ("field_license_wbddh") is filled in afterwards for all
# five rows, matching how it ended up last in the shared strings table.

$ws.Range("B270").Value = "CC0"
$ws.Range("C270").Value = "CC0 1.0"

$ws.Range("B271").Value = "Creative Commons Attribution 4.0"
$ws.Range("C271").Value = "Creative Commons Attribution 4.0"

$ws.Range("B272").Value = "Creative Commons Attribution-NonCommercial 4.0"
$ws.Range("C272").Value = "Creative Commons Attribution-NonCommercial 4.0"

$ws.Range("B273").Value = "Creative Commons Attribution Share-Alike 4.0"
$ws.Range("C273").Value = "Creative Commons Attribution Share-Alike 4.0"

$ws.Range("B274").Value = "Open Database License"
$ws.Range("C274").Value = "Open Data Commons Open Database License 1.0"

$ws.Range("A270").Value = "field_license_wbddh"
$ws.Range("A271").Value = "field_license_wbddh"
$ws.Range("A272").Value = "field_license_wbddh"
$ws.Range("A273").Value = "field_license_wbddh"
$ws.Range("A274").Value = "field_license_wbddh"

$ws.Range("D270").Value = $true
$ws.Range("D271").Value = $true
$ws.Range("D272").Value = $true
$ws.Range("D273").Value = $true
$ws.Range("D274").Value = $true

# Formatting for the new license-text cells: vertical-center + wrap text
# (matches the single new cellXfs entry added to styles.xml). Set the
# format on one cell first, then copy/paste-special just the formatting to
# the remaining cells so only one new style record is created instead of
# one per cell.
$firstFmt = $ws.Range("B270")
$firstFmt.WrapText = $true
$firstFmt.VerticalAlignment = -4108

$firstFmt.Copy() | Out-Null
$ws.Range("C270").PasteSpecial(-4122) | Out-Null
$ws.Range("B271:C274").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rows 272 and 274 hold longer license text that wraps onto a second line.
$ws.Rows.Item(272).RowHeight = 30
$ws.Rows.Item(274).RowHeight = 30

$ws.Range("I269").Select() | Out-Null
